# Append a new results row (row 7) that duplicates row 6's data
# ("Added half of the 24hr results, still feature selection left").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 6 (label + two accuracy numbers) down into row 7 so the new
# row reuses the same shared string / number values as the row above it.
$ws.Range("A6:C7").FillDown()

# FillDown also copies row 6's explicit column style onto row 7; nudge the
# font back to its existing value so Excel collapses the new cells back to
# the sheet's default (unstyled) cell format, matching the rest of the data
# rows (2-6), which likewise carry no explicit style.
$ws.Range("A7:C7").Font.Name = "Calibri"
